# "retraining scheme is working"
#
# Insert a new first data row into Table1 on Sheet1 for the new
# "LSTM w/ Retrain" portfolio (its Sharpe ratio isn't computed yet, so the
# Sharpe cell is left blank), shifting the existing LSTM/AGG/VTI/MVO/^VIX/DBC
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a blank row above the first data row (row 3); the table's
# DataBodyRange / old row 3.. shift down to row 4..
$ws.Rows(3).Insert()

# Pick up the number formatting of the (now-shifted) first data row so the
# new row matches the rest of the table (text column + 0.0000 numeric col).
$ws.Range("C4:D4").Copy()
$ws.Range("C3:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Grow the table definition to cover the new row.
$lo.Resize($ws.Range("C2:D9"))

# New portfolio name; Sharpe is left empty until the retrain run finishes.
$ws.Range("C3").Value = "LSTM w/ Retrain"

# Column C grew a bit wider to fit the new label.
$ws.Columns("C").ColumnWidth = 12.666666666666666

# Matches the selection left behind in the saved workbook.
[void]$ws.Range("C2:D9").Select()

# The saved file carries an explicit (portrait) page setup.
$ws.PageSetup.Orientation = 1
